$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnfsf14"
$ws.Cells.Item(2, 3).Value = "Tnfrsf14"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.02601766666666666
$ws.Cells.Item(2, 8).Value = 0.078053
$ws.Cells.Item(2, 9).Value = 0.02802616153292364
$ws.Cells.Item(2, 10).Value = 0.02802616153292364
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.667667333333333
$ws.Cells.Item(2, 14).Value = 11.003002
$ws.Cells.Item(2, 15).Value = 0.1424137080579054
$ws.Cells.Item(2, 16).Value = 0.1424137080579054
$ws.Cells.Item(2, 17).Value = 0.09542414612288887
$ws.Cells.Item(2, 18).Value = 0.8588173151059999
$ws.Cells.Item(2, 19).Value = 0.003991309586533485
$ws.Cells.Item(2, 20).Value = 0.003991309586533485

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnfsf14"
$ws.Cells.Item(3, 3).Value = "Tnfrsf14"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.02601766666666666
$ws.Cells.Item(3, 8).Value = 0.078053
$ws.Cells.Item(3, 9).Value = 0.02802616153292364
$ws.Cells.Item(3, 10).Value = 0.02802616153292364
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 6.497702
$ws.Cells.Item(3, 14).Value = 19.493106
$ws.Cells.Item(3, 15).Value = 0.2523025540689536
$ws.Cells.Item(3, 16).Value = 0.2523025540689536
$ws.Cells.Item(3, 17).Value = 0.1690550447353333
$ws.Cells.Item(3, 18).Value = 1.521495402618
$ws.Cells.Item(3, 19).Value = 0.007071072135505693
$ws.Cells.Item(3, 20).Value = 0.007071072135505694

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnfsf14"
$ws.Cells.Item(4, 3).Value = "Tnfrsf14"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.02601766666666666
$ws.Cells.Item(4, 8).Value = 0.078053
$ws.Cells.Item(4, 9).Value = 0.02802616153292364
$ws.Cells.Item(4, 10).Value = 0.02802616153292364
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 14.399313
$ws.Cells.Item(4, 14).Value = 43.197939
$ws.Cells.Item(4, 15).Value = 0.5591182000556945
$ws.Cells.Item(4, 16).Value = 0.5591182000556945
$ws.Cells.Item(4, 17).Value = 0.3746365258629999
$ws.Cells.Item(4, 18).Value = 3.371728732767
$ws.Cells.Item(4, 19).Value = 0.0156699369907584
$ws.Cells.Item(4, 20).Value = 0.01566993699075841

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Tnfsf14"
$ws.Cells.Item(5, 3).Value = "Tnfrsf14"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.02601766666666666
$ws.Cells.Item(5, 8).Value = 0.078053
$ws.Cells.Item(5, 9).Value = 0.02802616153292364
$ws.Cells.Item(5, 10).Value = 0.02802616153292364
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.188929333333334
$ws.Cells.Item(5, 14).Value = 3.566788
$ws.Cells.Item(5, 15).Value = 0.04616553781744659
$ws.Cells.Item(5, 16).Value = 0.04616553781744658
$ws.Cells.Item(5, 17).Value = 0.03093316708488889
$ws.Cells.Item(5, 18).Value = 0.278398503764
$ws.Cells.Item(5, 19).Value = 0.001293842820126053
$ws.Cells.Item(5, 20).Value = 0.001293842820126053

$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Tnfsf14"
$ws.Cells.Item(6, 3).Value = "Tnfrsf14"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.9023173333333333
$ws.Cells.Item(6, 8).Value = 2.706952
$ws.Cells.Item(6, 9).Value = 0.9719738384670763
$ws.Cells.Item(6, 10).Value = 0.9719738384670763
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.667667333333333
$ws.Cells.Item(6, 14).Value = 11.003002
$ws.Cells.Item(6, 15).Value = 0.1424137080579054
$ws.Cells.Item(6, 16).Value = 0.1424137080579054
$ws.Cells.Item(6, 17).Value = 3.309399807767111
$ws.Cells.Item(6, 18).Value = 29.784598269904
$ws.Cells.Item(6, 19).Value = 0.1384223984713719
$ws.Cells.Item(6, 20).Value = 0.1384223984713719

$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Tnfsf14"
$ws.Cells.Item(7, 3).Value = "Tnfrsf14"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.9023173333333333
$ws.Cells.Item(7, 8).Value = 2.706952
$ws.Cells.Item(7, 9).Value = 0.9719738384670763
$ws.Cells.Item(7, 10).Value = 0.9719738384670763
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 6.497702
$ws.Cells.Item(7, 14).Value = 19.493106
$ws.Cells.Item(7, 15).Value = 0.2523025540689536
$ws.Cells.Item(7, 16).Value = 0.2523025540689536
$ws.Cells.Item(7, 17).Value = 5.862989141434666
$ws.Cells.Item(7, 18).Value = 52.766902272912
$ws.Cells.Item(7, 19).Value = 0.2452314819334479
$ws.Cells.Item(7, 20).Value = 0.2452314819334479

$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Tnfsf14"
$ws.Cells.Item(8, 3).Value = "Tnfrsf14"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.9023173333333333
$ws.Cells.Item(8, 8).Value = 2.706952
$ws.Cells.Item(8, 9).Value = 0.9719738384670763
$ws.Cells.Item(8, 10).Value = 0.9719738384670763
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 14.399313
$ws.Cells.Item(8, 14).Value = 43.197939
$ws.Cells.Item(8, 15).Value = 0.5591182000556945
$ws.Cells.Item(8, 16).Value = 0.5591182000556945
$ws.Cells.Item(8, 17).Value = 12.992749707992
$ws.Cells.Item(8, 18).Value = 116.934747371928
$ws.Cells.Item(8, 19).Value = 0.543448263064936
$ws.Cells.Item(8, 20).Value = 0.543448263064936

$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Tnfsf14"
$ws.Cells.Item(9, 3).Value = "Tnfrsf14"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.9023173333333333
$ws.Cells.Item(9, 8).Value = 2.706952
$ws.Cells.Item(9, 9).Value = 0.9719738384670763
$ws.Cells.Item(9, 10).Value = 0.9719738384670763
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.188929333333334
$ws.Cells.Item(9, 14).Value = 3.566788
$ws.Cells.Item(9, 15).Value = 0.04616553781744659
$ws.Cells.Item(9, 16).Value = 0.04616553781744658
$ws.Cells.Item(9, 17).Value = 1.072791545575111
$ws.Cells.Item(9, 18).Value = 9.655123910176
$ws.Cells.Item(9, 19).Value = 0.04487169499732054
$ws.Cells.Item(9, 20).Value = 0.04487169499732053
